$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and 1h volume-change (E) figures.
# Price cells are forced to Text format before the write so that numeric-
# looking strings (e.g. "242.63") are not auto-converted to numbers by the
# Excel engine, matching the source data which stores them as plain text.
# The number format is reset back to Normal/General immediately afterwards
# so no residual cell-style change is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.938.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.952.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.43%  '

$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4876'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2942'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06957'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '107.11'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.976.02'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07755'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.50%  '

$ws.Range("E14").Value = '  -1.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6961'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '279.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.950.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007772'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.204.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.77%  '

$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("E22").Value = '  -2.72%  '

$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.485'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.733'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.172'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.93%  '

$ws.Range("E29").Value = '  -2.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.393'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.578'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.554'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.392'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04872'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7521'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.48%  '

$ws.Range("E36").Value = '  -0.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.723'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.80%  '

$ws.Range("E38").Value = '  -2.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.678'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.521'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '77.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.105'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8967'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4427'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9996'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.745'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '991.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1246'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.209'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.79'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.61%  '
